$d = $word.ActiveDocument
$d.Content.Find.Execute("Shift+Alt+1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Shift+Alt+0", 2)
